$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22; this shifts existing rows 22..67 down to 23..68
$ws.Rows.Item(22).Insert()

# Populate the newly-inserted row 22 with the new weekly record
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(22, 3).Value = "La Araucanía"
$ws.Cells.Item(22, 4).Value = 44972
$ws.Cells.Item(22, 5).Value = 9
$ws.Cells.Item(22, 6).Value = 100112042
$ws.Cells.Item(22, 7).Value = "Locoto"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 140
$ws.Cells.Item(22, 11).Value = 3300
$ws.Cells.Item(22, 12).Value = 3300
$ws.Cells.Item(22, 13).Value = 3300
$ws.Cells.Item(22, 14).Value = "$/kilo"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 3300
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
